# Updated symbol list on Sun Jan 15 11:53:26 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value (Price and Volume(1h) columns).
# A leading apostrophe forces Excel to store these as text (matching the
# original inlineStr/text cells) instead of auto-converting to numbers or
# percentages.
$updates = @{
    'D2' = "'296.73"
    'E2' = "'-2.08%"
    'D3' = "'31.27"
    'E3' = "'-1.45%"
    'D4' = "'5.110"
    'E4' = "'-1.84%"
    'D5' = "'0.07339"
    'E5' = "'-0.50%"
    'D6' = "'7.717"
    'E6' = "'-1.62%"
    'D7' = "'1.664"
    'E7' = "'10.31%"
    'D8' = "'3.728"
    'E8' = "'0.03%"
    'D9' = "'0.9178"
    'E9' = "'0.92%"
    'D10' = "'0.1665"
    'E10' = "'-0.66%"
    'D11' = "'0.07091"
    'E11' = "'-5.15%"
    'D12' = "'0.08072"
    'E12' = "'0.89%"
    'D13' = "'0.02987"
    'E13' = "'0.68%"
    'D14' = "'0.09900"
    'E14' = "'-0.11%"
    'D15' = "'0.001491"
    'E15' = "'-0.46%"
    'D16' = "'0.006182"
    'E16' = "'-1.05%"
    'E17' = "'-1.23%"
    'D18' = "'2.223"
    'E18' = "'-0.27%"
    'E20' = "'0.65%"
    'D21' = "'4.559"
    'E21' = "'1.67%"
    'D22' = "'0.04640"
    'E22' = "'2.40%"
    'D23' = "'0.1548"
    'E23' = "'-4.31%"
    'E24' = "'-0.04%"
    'D25' = "'0.004430"
    'E25' = "'0.16%"
    'D26' = "'0.0001299"
    'E26' = "'0.02%"
    'D27' = "'0.0001872"
    'E27' = "'7.66%"
    'D39' = "'0.01706"
    'E39' = "'2.22%"
    'D40' = "'0.04425"
    'E40' = "'-1.39%"
    'D41' = "'0.007215"
    'E41' = "'0.17%"
    'D42' = "'0.1330"
    'E42' = "'-0.98%"
    'D43' = "'0.002138"
    'E43' = "'-8.14%"
    'D44' = "'0.01112"
    'E44' = "'-22.25%"
    'D45' = "'0.00006006"
    'E45' = "'-2.18%"
    'E46' = "'-21.33%"
    'D47' = "'1.868"
    'E47' = "'-1.27%"
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

